$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '33.929.45'
$ws.Range("E2").Value = '  -0.61%  '
$ws.Range("D3").Value = '1.775.67'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.549'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.03'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.287'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0701'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0936'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = '2.034.58'
$ws.Range("E12").Value = '  -0.38%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.787.36'
$ws.Range("E13").Value = '  -0.33%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.87'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.77%  '
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '33.910.25'
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.618'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.13'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.72'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.54%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.39'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.84%  '
$ws.Range("D20").Value = '0.0₃0782'
$ws.Range("E20").Value = '  +0.49%  '
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.65'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.07'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.39%  '
$ws.Range("E24").Value = '  -2.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '160.34'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.25'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.112'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("E29").Value = '  +0.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.23'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0510'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.85%  '
$ws.Range("E32").Value = '  -1.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.80'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.85%  '
$ws.Range("D35").Value = '1.388.83'
$ws.Range("E35").Value = '  -1.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.652'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.91%  '
$ws.Range("E37").Value = '  -1.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0186'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.22'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.70%  '
$ws.Range("B40").Value = 'HuobiToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.36'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.907'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.91%  '
$ws.Range("E42").Value = '  -1.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '77.37'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.11'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +10.56%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.08'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.68%  '
$ws.Range("E46").Value = '  +12.02%  '
$ws.Range("B47").Value = 'Kaspa'
$ws.Range("C47").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0495'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.80%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '107.71'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.15%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.81'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.54%  '
$ws.Range("D50").Value = '1.932.08'
$ws.Range("E50").Value = '  -0.37%  '
$ws.Range("E51").Value = '  +0.66%  '
